$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "MMR_inventory_template"

# Remove the now-unused "type"/"smr" column (G)
$ws.Columns.Item(7).Delete()

# Update header labels (order matters for shared-string table layout:
# "start"/"end" must be registered before "cycle")
$ws.Range("E1").Value = "start"
$ws.Range("F1").Value = "end"
$ws.Range("D1").Value = "cycle"

# Row 2 data
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 10.9
$ws.Range("F2").Value = 13

# Row 3 data
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 10.3
$ws.Range("F3").Value = 14

# Row 4: clear everything, leave A4 as an empty, quote-prefixed cell
$ws.Range("A4").Style = "Normal"
$ws.Range("A4").Value = "'x"
$ws.Range("A4").ClearContents()
$ws.Range("B4:F4").ClearContents()

# Drop the old row 5 (no longer needed)
$ws.Rows.Item(5).Delete()

# Widen columns E:L to match the new "start"/"end" data columns
$ws.Range($ws.Cells.Item(1,5), $ws.Cells.Item(1,12)).EntireColumn.ColumnWidth = 15.8333333

# Update the selection shown when the sheet is reopened
$null = $ws.Range("E6").Select()

Write-Host "done"
